# Inventario.xlsx update
# - Adds a new worksheet ("Sheet1") after "Java Books" with a small
#   A1:C3 block of single-letter test values.
# - Appends 16 more book rows (No. 11-26) to the "Java Books" sheet,
#   cycling through 4 new book/author/price combinations 4 times.

$wb = $excel.ActiveWorkbook

# --- 1) Add the new "Sheet1" worksheet after the existing last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$sheet2.Name = "Sheet1"

# Fill it in row-major order so shared strings are interned in the
# same order they are first seen (q, l, h, o, a).
$sheet2.Range("A1").Value = "q"
$sheet2.Range("B1").Value = "l"
$sheet2.Range("C1").Value = "q"
$sheet2.Range("A2").Value = "h"
$sheet2.Range("B2").Value = "o"
$sheet2.Range("C2").Value = "l"
$sheet2.Range("A3").Value = "a"
$sheet2.Range("B3").Value = "a"
$sheet2.Range("C3").Value = "a"
$sheet2.Range("C3").Select()

# --- 2) Append the new book rows to "Java Books" ---
$ws = $wb.Worksheets.Item("Java Books")

$books = @(
    @("El que se duerme pierde", "Tom Peter", 16),
    @("Sin lugar a duda", "Ana Gutierrez", 26),
    @("El arte de dormir", "Nico", 32),
    @("Buscando a Nemo", "Humble Po", 41)
)

$row = 12
$no = 11
for ($rep = 0; $rep -lt 4; $rep++) {
    foreach ($b in $books) {
        $ws.Cells.Item($row, 1).Value = $no
        $ws.Cells.Item($row, 2).Value = $b[0]
        $ws.Cells.Item($row, 3).Value = $b[1]
        $ws.Cells.Item($row, 4).Value = $b[2]
        $row++
        $no++
    }
}

# Leave the original sheet active, with the selection where the
# author's session ended up.
$ws.Activate()
$ws.Range("I19").Select()
